$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New work-log entries. Row 23 (date/hours) already had its "Arbeit" text (G23)
# filled in, but the Datum / Zeit in h cells were still empty; rows 24-32 are
# brand new entries appended below it.

$rows = @(
    @{ Row = 23; Date = 43685; Hours = 0.5;  Text = "Funktionalität der UI programmiert" },
    @{ Row = 24; Date = 43689; Hours = 2.5;  Text = "Funktionalität der UI programmiert, Probleme beim Zuweisen von Eigenschaften im Scene Builder" },
    @{ Row = 25; Date = 43690; Hours = 1;    Text = "Statusupdate erstellt" },
    @{ Row = 26; Date = 43692; Hours = 0.25; Text = "Fehlersuche in Scene Builder" },
    @{ Row = 27; Date = 43695; Hours = 0.5;  Text = "Fehlersuche in Scene Builder" },
    @{ Row = 28; Date = 43699; Hours = 0.2;  Text = "Aufbau der UI überdacht, neuer Menüpunkt eingefügt" },
    @{ Row = 29; Date = 43700; Hours = 0.25; Text = "Neuer Menüpunkt eingefügt" },
    @{ Row = 30; Date = 43705; Hours = 0.5;  Text = "Statusupdate erstellt" },
    @{ Row = 31; Date = 43711; Hours = 1.5;  Text = "Funktionalität der UI programmiert" },
    @{ Row = 32; Date = 43712; Hours = 2.5;  Text = "Funktionalität der UI programmiert" }
)

# Row 21 carries the "normal" (non-wrapped, non-tall) formatting for the
# E:G block -- use it as the template for every newly appended row so the
# copied styles match the existing sheet instead of creating fresh ones.
# Values are written first (so the new row/cells exist before any
# formatting copy touches them) and the formats are copied over afterward;
# this keeps the SUM(F:F) dependency on the freshly-created cells intact.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Date
    $ws.Cells.Item($r.Row, 6).Value = $r.Hours
    $ws.Cells.Item($r.Row, 7).Value = $r.Text

    $ws.Range("E21:F21").Copy()
    $ws.Range("E" + $r.Row + ":F" + $r.Row).PasteSpecial(-4122)
    $ws.Range("G21").Copy()
    $ws.Range("G" + $r.Row).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update sheet view scroll position / selection to match the edited area.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("G29").Select()

$excel.Calculate()
